# Scheduled market-price refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit columns (H, I, J, K, L, M, N) across the per-job
# sheets, row by row, to reflect freshly pulled Market Board data.
# Where HQ pricing data is no longer available for a leve, the HQ profit
# cell (column N, occasionally M) is cleared entirely rather than zeroed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1236.75
$ws.Range("J19").Value = 724.2
$ws.Range("L19").Value = 724.2
$ws.Range("N19").Value = -1074.2

$ws.Range("H43").Value = 4429.067
$ws.Range("I43").Value = 4118.3335
$ws.Range("J43").Value = 5154.1113
$ws.Range("K43").Value = 4118.3335
$ws.Range("L43").Value = 5154.1113
$ws.Range("M43").Value = -4049.3335
$ws.Range("N43").Value = -5292.1113

$ws.Range("H64").Value = 7198
$ws.Range("I64").Value = 3463.8333
$ws.Range("K64").Value = 3463.8333
$ws.Range("M64").Value = -3215.8333

$ws.Range("H67").Value = 7198
$ws.Range("I67").Value = 3463.8333
$ws.Range("K67").Value = 3463.8333
$ws.Range("M67").Value = -2605.8333

$ws.Range("H74").Value = 2449.3333
$ws.Range("I74").Value = 2449.3333
$ws.Range("K74").Value = 2449.3333
$ws.Range("M74").Value = -1513.3333

$ws.Range("H77").Value = 2449.3333
$ws.Range("I77").Value = 2449.3333
$ws.Range("K77").Value = 12246.6665
$ws.Range("M77").Value = -7566.666499999999

$ws.Range("H80").Value = 1376.1482
$ws.Range("I80").Value = 906.4
$ws.Range("J80").Value = 1652.4706
$ws.Range("K80").Value = 2719.2
$ws.Range("L80").Value = 4957.4118
$ws.Range("M80").Value = -1721.2
$ws.Range("N80").Value = -6953.4118

$ws.Range("H83").Value = 1376.1482
$ws.Range("I83").Value = 906.4
$ws.Range("J83").Value = 1652.4706
$ws.Range("K83").Value = 8157.599999999999
$ws.Range("L83").Value = 14872.2354
$ws.Range("M83").Value = -3165.599999999999
$ws.Range("N83").Value = -24856.2354

$ws.Range("H86").Value = 7893.2856
$ws.Range("I86").Value = 6900
$ws.Range("J86").Value = 8638.25
$ws.Range("K86").Value = 6900
$ws.Range("L86").Value = 8638.25
$ws.Range("M86").Value = -5777
$ws.Range("N86").Value = -10884.25

$ws.Range("H89").Value = 7893.2856
$ws.Range("I89").Value = 6900
$ws.Range("J89").Value = 8638.25
$ws.Range("K89").Value = 34500
$ws.Range("L89").Value = 43191.25
$ws.Range("M89").Value = -28884
$ws.Range("N89").Value = -54423.25

$ws.Range("H113").Value = 16455
$ws.Range("J113").Value = 10002.5
$ws.Range("L113").Value = 10002.5
$ws.Range("N113").Value = -16510.5

$ws.Range("H132").Value = 1776.6216
$ws.Range("I132").Value = 1460.7812
$ws.Range("J132").Value = 3798
$ws.Range("K132").Value = 4382.3436
$ws.Range("L132").Value = 11394
$ws.Range("M132").Value = -1852.3436
$ws.Range("N132").Value = -16454

$ws.Range("H137").Value = 1363.6086
$ws.Range("I137").Value = 1254.4706
$ws.Range("J137").Value = 1672.8334
$ws.Range("K137").Value = 3763.4118
$ws.Range("L137").Value = 5018.5002
$ws.Range("M137").Value = -1213.4118
$ws.Range("N137").Value = -10118.5002

$ws.Range("H138").Value = 3520
$ws.Range("I138").Value = 3520
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 10560
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -5420
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3844
$ws.Range("I122").Value = 2399.5
$ws.Range("J122").Value = 5288.5
$ws.Range("K122").Value = 7198.5
$ws.Range("L122").Value = 15865.5
$ws.Range("M122").Value = -4748.5
$ws.Range("N122").Value = -20765.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2751.3845
$ws.Range("I86").Value = 1935.7142
$ws.Range("K86").Value = 1935.7142
$ws.Range("M86").Value = -812.7141999999999

$ws.Range("H89").Value = 2751.3845
$ws.Range("I89").Value = 1935.7142
$ws.Range("K89").Value = 9678.571
$ws.Range("M89").Value = -4062.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 13214145
$ws.Range("I86").Value = 13214145
$ws.Range("K86").Value = 13214145
$ws.Range("M86").Value = -13213022

$ws.Range("H89").Value = 13214145
$ws.Range("I89").Value = 13214145
$ws.Range("K89").Value = 66070725
$ws.Range("M89").Value = -66065109

$ws.Range("H99").Value = 5942
$ws.Range("I99").Value = 6124.6
$ws.Range("K99").Value = 6124.6
$ws.Range("M99").Value = -4626.6

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H126").Value = 5942
$ws.Range("I126").Value = 6124.6
$ws.Range("K126").Value = 18373.8
$ws.Range("M126").Value = -15903.8

$ws.Range("H134").Value = 1826.6666
$ws.Range("I134").Value = 1568.1538
$ws.Range("K134").Value = 4704.4614
$ws.Range("M134").Value = -2169.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 995
$ws.Range("I133").Value = 995
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 2985
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = 2075
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11681.8
$ws.Range("I70").Value = 11750
$ws.Range("K70").Value = 11750
$ws.Range("M70").Value = -11480

$ws.Range("H73").Value = 11681.8
$ws.Range("I73").Value = 11750
$ws.Range("K73").Value = 11750
$ws.Range("M73").Value = -10814

$ws.Range("H80").Value = 10500
$ws.Range("J80").Value = 11752.5
$ws.Range("L80").Value = 11752.5
$ws.Range("N80").Value = -13748.5

$ws.Range("H83").Value = 10500
$ws.Range("J83").Value = 11752.5
$ws.Range("L83").Value = 58762.5
$ws.Range("N83").Value = -68746.5

$ws.Range("H122").Value = 4622.25
$ws.Range("J122").Value = 4000
$ws.Range("L122").Value = 12000
$ws.Range("N122").Value = -16900

$ws.Range("H123").Value = 33499.5
$ws.Range("J123").Value = 33499.5
$ws.Range("L123").Value = 33499.5
$ws.Range("N123").Value = -38399.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3251
$ws.Range("I40").Value = 3251
$ws.Range("K40").Value = 3251
$ws.Range("M40").Value = -3115

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3374.75
$ws.Range("I122").Value = 3374.75
$ws.Range("K122").Value = 10124.25
$ws.Range("M122").Value = -7674.25

$ws.Range("H136").Value = 754.0833
$ws.Range("I136").Value = 754.0833
$ws.Range("K136").Value = 2262.2499
$ws.Range("M136").Value = 287.7501000000002
